# Rename the inline logo pictures in the document's headers/footers.
#
#   headers (BTec_Logo-Orange, .jpg):  image2.jpg -> image1.jpg
#   footers (PearsonLogo, .png):       image1.png -> image2.png
#
# InlineShape has no high-fidelity "Name" property of its own via the
# Range path once the handle has gone through a structural re-fetch (this
# shows up for footers in particular), so the shape is first Select()-ed
# and then addressed again through Selection.InlineShapes before the
# rename is applied. That sequence reliably lands the write.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Headers: BTec logo, image2.jpg -> image1.jpg -------------------------
$headers = $sec.Headers
for ($i = 1; $i -le $headers.Count; $i++) {
    $hf = $headers.Item($i)
    if ($hf.Exists) {
        $shapes = $hf.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shapes.Item($j).Name = "image1.jpg"
        }
    }
}

# --- Footers: Pearson logo, image1.png -> image2.png -----------------------
$footers = $sec.Footers
for ($i = 1; $i -le $footers.Count; $i++) {
    $hf = $footers.Item($i)
    if ($hf.Exists) {
        $shapes = $hf.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shape = $shapes.Item($j)
            $shape.Select() | Out-Null
            $word.Selection.InlineShapes.Item(1).Name = "image2.png"
        }
    }
}

Write-Output "Renamed header/footer logo pictures."
